# ---------------------------------------------------------------------------
# Live trading results update
#
#  - Trade #51 (leadlag, opened 2026-02-16 21:31:02) is closed out: exit
#    price, P&L, exit reason and duration are filled in on the "leadlag"
#    sheet (row 41) and the now-closed trade is appended as a new row on
#    "All Trades" (row 52).
#  - A new trade #75 (momentum, opened 2026-02-16 21:36:03) is appended as
#    a freshly OPENed row on "momentum" (row 20).
#  - Summary / Comparison aggregate stats are recalculated accordingly.
#
# Several of the text values below ("66.7%", "2.80", "+0.5573%", ...) look
# like numbers/percentages, and Excel's normal cell-typing heuristics would
# silently reinterpret them as numeric values. The source data keeps them
# as literal text, so the ranges that need it are pre-formatted as Text
# ("@") before the value is typed in, which keeps the literal string
# (matching how the workbook already stores its other text-shaped numbers).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# Pre-format the percentage-shaped text cells as Text so typing the value
# doesn't get auto-converted into a numeric percentage.
$wsSummary.Range("D2:F3").NumberFormat = "@"

# Row 2: OVERALL / ALL COMBINED
$wsSummary.Range("C2").Value = 51
$wsSummary.Range("D2").Value = "66.7%"
$wsSummary.Range("E2").Value = "+13.1609%"
$wsSummary.Range("F2").Value = "+0.2581%"

# Row 3: STRATEGY / leadlag
$wsSummary.Range("C3").Value = 56
$wsSummary.Range("D3").Value = "44.6%"
$wsSummary.Range("E3").Value = "+8.9593%"
$wsSummary.Range("F3").Value = "+0.1600%"

# ---------------------------------------------------------------------------
# leadlag sheet - close out trade #51 (row 41)
# ---------------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

$wsLeadlag.Range("G41").Value = 68335.91757400001
$wsLeadlag.Range("H41").Value = "CLOSED"
$wsLeadlag.Range("I41").Value = 0.6082
$wsLeadlag.Range("J41").Value = 6.08
$wsLeadlag.Range("M41").Value = "time_exit_5min"
$wsLeadlag.Range("N41").Value = 5

# ---------------------------------------------------------------------------
# momentum sheet - append newly opened trade #75 (row 20)
# ---------------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

# "2026-02-16" reads as a date to Excel's auto-detection, so keep it literal
# text the same way.
$wsMomentum.Range("B20").NumberFormat = "@"

$wsMomentum.Range("A20").Value = 75
$wsMomentum.Range("B20").Value = "2026-02-16"
$wsMomentum.Range("C20").Value = "21:36:03"
$wsMomentum.Range("D20").Value = "momentum"
$wsMomentum.Range("E20").Value = "UP"
$wsMomentum.Range("F20").Value = 68700.565
$wsMomentum.Range("H20").Value = "OPEN"
$wsMomentum.Range("I20").Value = 0
$wsMomentum.Range("J20").Value = 0
$wsMomentum.Range("K20").Value = 0.9
$wsMomentum.Range("L20").Value = "Upward momentum: 0.317% over 10 samples"
$wsMomentum.Range("N20").Value = 0
# G20 / M20 stay blank (trade is still open, no exit info yet).

# ---------------------------------------------------------------------------
# All Trades sheet - append the now-closed trade #51 (row 52)
# ---------------------------------------------------------------------------
$wsAllTrades = $wb.Worksheets.Item("All Trades")

$wsAllTrades.Range("B52").NumberFormat = "@"

$wsAllTrades.Range("A52").Value = 51
$wsAllTrades.Range("B52").Value = "2026-02-16"
$wsAllTrades.Range("C52").Value = "21:31:02"
$wsAllTrades.Range("D52").Value = "leadlag"
$wsAllTrades.Range("E52").Value = "DOWN"
$wsAllTrades.Range("F52").Value = 68754.095
$wsAllTrades.Range("G52").Value = 68335.91757400001
$wsAllTrades.Range("H52").Value = "CLOSED"
$wsAllTrades.Range("I52").Value = 0.6082
$wsAllTrades.Range("J52").Value = 6.08
$wsAllTrades.Range("K52").Value = 0.7212
$wsAllTrades.Range("L52").Value = "Coinbase leading with -0.072% move"
$wsAllTrades.Range("M52").Value = "time_exit_5min"
$wsAllTrades.Range("N52").Value = 5

# ---------------------------------------------------------------------------
# Comparison sheet
# ---------------------------------------------------------------------------
$wsComparison = $wb.Worksheets.Item("Comparison")

$wsComparison.Range("C2").NumberFormat = "@"
$wsComparison.Range("D2").NumberFormat = "@"
$wsComparison.Range("E2").NumberFormat = "@"
$wsComparison.Range("G2").NumberFormat = "@"

$wsComparison.Range("B2").Value = 56
$wsComparison.Range("C2").Value = "44.6%"
$wsComparison.Range("D2").Value = "2.80"
$wsComparison.Range("E2").Value = "+0.5573%"
$wsComparison.Range("G2").Value = "1.68"
